# Updates the crypto price table (columns D = Price, E = Volume(1h))
# for rows 2-51 of the active sheet, matching the source data refresh.
#
# Both columns hold plain TEXT in the workbook (not numbers/percentages),
# so every new value must land back in the cell as a string. A bare
# `.Value = '123.45'` assignment is unsafe for plain numeric-looking
# strings: Excel's COM layer auto-coerces those into real numbers
# (dropping e.g. trailing zeros), which would silently change the cell
# type away from text. Column E is never at risk (its values carry
# padding spaces / a trailing '%', so Excel leaves them as text), and a
# handful of column D values are safe too (they contain two '.' like
# '26.165.89', which can't parse as a number either). For the remaining,
# genuinely numeric-looking D values we force text by building the exact
# decimal string with TEXT(), writing it through an off-sheet scratch
# cell, and landing only the computed VALUE via PasteSpecial(xlPasteValues)
# -- this avoids leaving any NumberFormat/style behind on the target cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Off-sheet scratch cell (outside the A1:E51 data range) used as a relay
# to push a TEXT()-formula result into a target cell as a literal string.
$scratch = $ws.Cells.Item(100, 26)

function Set-TextValue($cell, $value) {
    $cell.Value = $value
}

function Set-NumericLookingTextValue($cell, $numericText, $decimalFormat) {
    # Round-trip the literal through TEXT()+PasteSpecial so Excel can't
    # reinterpret the numeric-looking string as a real number.
    $scratch.Formula = '=TEXT(' + $numericText + ',"' + $decimalFormat + '")'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

Set-TextValue $ws.Cells.Item(2, 4) '26.165.89'
Set-TextValue $ws.Cells.Item(2, 5) '  -2.12%  '
Set-TextValue $ws.Cells.Item(3, 4) '1.668.44'
Set-TextValue $ws.Cells.Item(3, 5) '  -1.57%  '
Set-TextValue $ws.Cells.Item(4, 5) '  -0.08%  '
Set-NumericLookingTextValue $ws.Cells.Item(5, 4) '216.92' '0.00'
Set-TextValue $ws.Cells.Item(5, 5) '  -0.62%  '
Set-TextValue $ws.Cells.Item(6, 5) '  +2.28%  '
Set-NumericLookingTextValue $ws.Cells.Item(7, 4) '1.006' '0.000'
Set-TextValue $ws.Cells.Item(7, 5) '  -0.03%  '
Set-NumericLookingTextValue $ws.Cells.Item(8, 4) '0.2639' '0.0000'
Set-TextValue $ws.Cells.Item(8, 5) '  +1.91%  '
Set-TextValue $ws.Cells.Item(9, 5) '  +4.90%  '
Set-TextValue $ws.Cells.Item(10, 5) '  -0.58%  '
Set-NumericLookingTextValue $ws.Cells.Item(11, 4) '0.07408' '0.00000'
Set-TextValue $ws.Cells.Item(11, 5) '  +1.42%  '
Set-TextValue $ws.Cells.Item(12, 4) '1.671.16'
Set-TextValue $ws.Cells.Item(12, 5) '  -1.22%  '
Set-NumericLookingTextValue $ws.Cells.Item(13, 4) '4.516' '0.000'
Set-TextValue $ws.Cells.Item(13, 5) '  +2.11%  '
Set-NumericLookingTextValue $ws.Cells.Item(14, 4) '0.5816' '0.0000'
Set-TextValue $ws.Cells.Item(14, 5) '  +1.83%  '
Set-NumericLookingTextValue $ws.Cells.Item(15, 4) '0.000008564' '0.000000000'
Set-TextValue $ws.Cells.Item(15, 5) '  +4.77%  '
Set-NumericLookingTextValue $ws.Cells.Item(16, 4) '64.25' '0.00'
Set-TextValue $ws.Cells.Item(16, 5) '  -0.99%  '
Set-TextValue $ws.Cells.Item(17, 4) '26.227.92'
Set-TextValue $ws.Cells.Item(17, 5) '  -1.96%  '
Set-NumericLookingTextValue $ws.Cells.Item(18, 4) '4.939' '0.000'
Set-TextValue $ws.Cells.Item(18, 5) '  -1.42%  '
Set-TextValue $ws.Cells.Item(19, 5) '  -0.06%  '
Set-TextValue $ws.Cells.Item(20, 5) '  +0.90%  '
Set-NumericLookingTextValue $ws.Cells.Item(21, 4) '190.49' '0.00'
Set-TextValue $ws.Cells.Item(21, 5) '  +3.69%  '
Set-NumericLookingTextValue $ws.Cells.Item(22, 4) '6.214' '0.000'
Set-TextValue $ws.Cells.Item(22, 5) '  +0.22%  '
Set-NumericLookingTextValue $ws.Cells.Item(23, 4) '1.007' '0.000'
Set-TextValue $ws.Cells.Item(23, 5) '  +0.03%  '
Set-NumericLookingTextValue $ws.Cells.Item(24, 4) '145.55' '0.00'
Set-TextValue $ws.Cells.Item(24, 5) '  +0.20%  '
Set-NumericLookingTextValue $ws.Cells.Item(25, 4) '7.630' '0.000'
Set-TextValue $ws.Cells.Item(25, 5) '  +0.36%  '
Set-NumericLookingTextValue $ws.Cells.Item(26, 4) '0.1196' '0.0000'
Set-TextValue $ws.Cells.Item(26, 5) '  +5.19%  '
Set-NumericLookingTextValue $ws.Cells.Item(27, 4) '15.62' '0.00'
Set-TextValue $ws.Cells.Item(27, 5) '  +2.51%  '
Set-NumericLookingTextValue $ws.Cells.Item(28, 4) '0.06370' '0.00000'
Set-TextValue $ws.Cells.Item(28, 5) '  +13.93%  '
Set-NumericLookingTextValue $ws.Cells.Item(29, 4) '1.296' '0.000'
Set-TextValue $ws.Cells.Item(29, 5) '  -1.39%  '
Set-NumericLookingTextValue $ws.Cells.Item(30, 4) '1.319' '0.000'
Set-TextValue $ws.Cells.Item(30, 5) '  -0.94%  '
Set-TextValue $ws.Cells.Item(31, 5) '  +2.25%  '
Set-NumericLookingTextValue $ws.Cells.Item(32, 4) '3.522' '0.000'
Set-TextValue $ws.Cells.Item(32, 5) '  +2.34%  '
Set-NumericLookingTextValue $ws.Cells.Item(33, 4) '1.638' '0.000'
Set-TextValue $ws.Cells.Item(33, 5) '  -0.93%  '
Set-NumericLookingTextValue $ws.Cells.Item(34, 4) '1.014' '0.000'
Set-TextValue $ws.Cells.Item(34, 5) '  +0.83%  '
Set-NumericLookingTextValue $ws.Cells.Item(35, 4) '0.6085' '0.0000'
Set-TextValue $ws.Cells.Item(35, 5) '  +3.65%  '
Set-NumericLookingTextValue $ws.Cells.Item(36, 4) '2.368' '0.000'
Set-TextValue $ws.Cells.Item(36, 5) '  -1.58%  '
Set-NumericLookingTextValue $ws.Cells.Item(37, 4) '2.657' '0.000'
Set-TextValue $ws.Cells.Item(37, 5) '  +1.01%  '
Set-TextValue $ws.Cells.Item(38, 5) '  +4.85%  '
Set-TextValue $ws.Cells.Item(39, 5) '  +1.55%  '
Set-TextValue $ws.Cells.Item(40, 4) '1.083.26'
Set-TextValue $ws.Cells.Item(40, 5) '  +1.32%  '
Set-NumericLookingTextValue $ws.Cells.Item(41, 4) '0.8665' '0.0000'
Set-TextValue $ws.Cells.Item(41, 5) '  +1.82%  '
Set-NumericLookingTextValue $ws.Cells.Item(43, 4) '101.38' '0.00'
Set-TextValue $ws.Cells.Item(43, 5) '  +2.86%  '
Set-TextValue $ws.Cells.Item(44, 4) '1.817.37'
Set-TextValue $ws.Cells.Item(44, 5) '  -1.99%  '
Set-TextValue $ws.Cells.Item(45, 5) '  +5.23%  '
Set-NumericLookingTextValue $ws.Cells.Item(46, 4) '56.23' '0.00'
Set-TextValue $ws.Cells.Item(46, 5) '  -0.04%  '
Set-TextValue $ws.Cells.Item(47, 5) '  +0.44%  '
Set-NumericLookingTextValue $ws.Cells.Item(48, 4) '8.124' '0.000'
Set-TextValue $ws.Cells.Item(48, 5) '  +0.17%  '
Set-NumericLookingTextValue $ws.Cells.Item(49, 4) '0.05204' '0.00000'
Set-TextValue $ws.Cells.Item(49, 5) '  -0.05%  '
Set-NumericLookingTextValue $ws.Cells.Item(50, 4) '0.4294' '0.0000'
Set-TextValue $ws.Cells.Item(50, 5) '  -0.95%  '
Set-NumericLookingTextValue $ws.Cells.Item(51, 4) '5.898' '0.000'
Set-TextValue $ws.Cells.Item(51, 5) '  +5.52%  '

$scratch.ClearContents()
